$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a note/footnote below the country table explaining the assumption
# methodology (row 37 stays empty, content lands in row 38).
$ws.Range("A38").Value = "Note:"
$ws.Range("B38").Value = "If there is an assumption, the energy per capita of the compensatory country was multiplied with population of the missing country"

# Leave the selection where the user ended up after typing the note.
$ws.Range("B39").Select()
